# Applies updated First_Noticeable_Increase_Index (C), First_Noticeable_Increase_Cumulative_Value (E)
# and Pulse_Width (G) values on each Step3_DataPts_* sheet, reflecting the newly introduced
# zero_before_threshold parameter (values before the noise threshold / first rise point are
# now zeroed out before the index/cumulative-value/pulse-width calculations).

$wb = $excel.ActiveWorkbook

# Row -> (C: First_Noticeable_Increase_Index, E: First_Noticeable_Increase_Cumulative_Value, G: Pulse_Width)
# These are identical for columns C/E across all four threshold sheets; G (Pulse_Width = D - C)
# varies per sheet because D (Point_Exceeds_Index) differs per Intensity_Threshold.
$updates = @{
    "Step3_DataPts_0.5" = @{
        2 = @(89, 0.05006871881067964, 20)
        3 = @(88, 0.01531230265387533, 43)
        4 = @(90, 0.01252093636086236, 41)
        5 = @(89, 0.03469145590822811, 24)
        6 = @(88, 0.02897671090433112, 38)
    }
    "Step3_DataPts_0.7" = @{
        2 = @(89, 0.05006871881067964, 55)
        3 = @(88, 0.01531230265387533, 61)
        4 = @(90, 0.01252093636086236, 59)
        5 = @(89, 0.03469145590822811, 58)
        6 = @(88, 0.02897671090433112, 62)
    }
    "Step3_DataPts_0.8" = @{
        2 = @(89, 0.05006871881067964, 68)
        3 = @(88, 0.01531230265387533, 80)
        4 = @(90, 0.01252093636086236, 78)
        5 = @(89, 0.03469145590822811, 69)
        6 = @(88, 0.02897671090433112, 87)
    }
    "Step3_DataPts_0.9" = @{
        2 = @(89, 0.05006871881067964, 87)
        3 = @(88, 0.01531230265387533, 106)
        4 = @(90, 0.01252093636086236, 104)
        5 = @(89, 0.03469145590822811, 98)
        6 = @(88, 0.02897671090433112, 107)
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $vals = $rows[$rowNum]
        $ws.Cells.Item($rowNum, 3).Value = $vals[0]   # Column C
        $ws.Cells.Item($rowNum, 5).Value = $vals[1]   # Column E
        $ws.Cells.Item($rowNum, 7).Value = $vals[2]   # Column G
    }
}

Write-Host "Updated zero_before_threshold derived values on Step3_DataPts sheets."
